$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.378.18"
$ws.Range("E2").Value = "  -0.39%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.944.34"
$ws.Range("E3").Value = "  -2.06%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.19"
$ws.Range("E5").Value = "  -0.87%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.611"
$ws.Range("E6").Value = "  -1.55%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.97"
$ws.Range("E8").Value = "  -3.78%  "

$ws.Range("E9").Value = "  -4.02%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0852"
$ws.Range("E10").Value = "  +4.07%  "

$ws.Range("E11").Value = "  +0.27%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.230.11"
$ws.Range("E12").Value = "  -2.00%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.815"
$ws.Range("E13").Value = "  -5.56%  "

$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "13.47"
$ws.Range("E14").Value = "  -4.21%  "

$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.07"
$ws.Range("E15").Value = "  -12.06%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.16"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.961.88"
$ws.Range("E17").Value = "  -1.55%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.342.16"
$ws.Range("E18").Value = "  -0.23%  "

$ws.Range("E19").Value = "  +1.00%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.37"
$ws.Range("E20").Value = "  -1.92%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "228.74"
$ws.Range("E21").Value = "  -2.09%  "

$ws.Range("E22").Value = "  -6.16%  "

$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.38"
$ws.Range("E24").Value = "  -8.79%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.27"
$ws.Range("E25").Value = "  -1.54%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.22"
$ws.Range("E26").Value = "  -10.11%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.37"
$ws.Range("E27").Value = "  -0.36%  "

$ws.Range("E28").Value = "  +4.28%  "

$ws.Range("E29").Value = "  -3.46%  "

$ws.Range("E30").Value = "  -2.28%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.13"
$ws.Range("E31").Value = "  -5.80%  "

$ws.Range("E32").Value = "  -6.18%  "

$ws.Range("E33").Value = "  +0.53%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.24"
$ws.Range("E34").Value = "  -4.60%  "

$ws.Range("B35").Value = "BinanceUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.03%  "

$ws.Range("B36").Value = "THORChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.11"
$ws.Range("E36").Value = "  -3.47%  "

$ws.Range("E37").Value = "  +1.32%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.14"
$ws.Range("E38").Value = "  -6.01%  "

$ws.Range("E39").Value = "  -2.87%  "

$ws.Range("E40").Value = "  +0.09%  "

$ws.Range("E41").Value = "  -1.18%  "

$ws.Range("E42").Value = "  -1.73%  "

$ws.Range("E43").Value = "  -7.23%  "

$ws.Range("E44").Value = "  -4.06%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.344.87"
$ws.Range("E45").Value = "  -2.57%  "

$ws.Range("E46").Value = "  -6.81%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "87.08"
$ws.Range("E47").Value = "  -6.47%  "

$ws.Range("E48").Value = "  -6.38%  "

$ws.Range("E49").Value = "  -1.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "44.27"
$ws.Range("E50").Value = "  -1.80%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.120.57"
$ws.Range("E51").Value = "  -2.29%  "
